{"js": "// Update the year in the astromap credit-line URL from 2018 to 2022.\n// The credit paragraph (\"Les cartes figurant dans ce document ont \u00e9t\u00e9\n// \u00e9tablies par Jenik Hollan, CzechGlobe ((http://.../GaNight/2018/).\")\n// is split across many small runs (one per word/character, with\n// proofing marks in between). Find that paragraph, then replace its\n// whole range with the corrected text, collapsing it to a single run.\n\nconst body = context.document.body;\nconst results = body.search(\"Les cartes figurant dans ce document\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the astromap credit paragraph.\");\n}\n\nconst paragraph = results.items[0].paragraphs.getFirst();\nconst paragraphRange = paragraph.getRange();\n\nconst newText =\n  \"Les cartes figurant dans ce document ont \u00e9t\u00e9 \u00e9tablies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\";\n\n// Clear existing run formatting first so the replacement text lands in\n// a single plain run (matching how the source edit collapsed the runs).\nparagraphRange.clear();\nawait context.sync();\n\nparagraphRange.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the year in the astromap credit-line URL from 2018 to 2022.\n# The credit paragraph (\"Les cartes figurant dans ce document ont \u00e9t\u00e9\n# \u00e9tablies par Jenik Hollan, CzechGlobe ((http://.../GaNight/2018/).\")\n# is made up of many small runs (one per word/character, with proofing\n# marks in between). We replace the whole paragraph's text (but not its\n# paragraph mark, so paragraph formatting / borders / shading are kept)\n# with the same text but with the year corrected, collapsing it down to\n# a single run.\n\n$d = $word.ActiveDocument\n\n$oldText = \"Les cartes figurant dans ce document ont \u00e9t\u00e9 \u00e9tablies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).\"\n$newText = \"Les cartes figurant dans ce document ont \u00e9t\u00e9 \u00e9tablies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\"\n\nforeach ($p in $d.Paragraphs) {\n    $rng = $p.Range\n    if ($rng.Text -like \"*Les cartes figurant dans ce document*\") {\n        # Shrink the range so it excludes the trailing paragraph mark,\n        # then swap its contents for the corrected credit line.\n        $rng.End = $rng.End - 1\n        $rng.Delete()\n        $rng.InsertAfter($newText)\n        break\n    }\n}\n"}
